# Update column E (rows 2 to 23) from 50 to 70 on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2:E23").Value = 70
